$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency list (coin, link, price, 1h volume) for rows 2-51.
# A new coin (WrappedeETH) was inserted at row 28, shifting subsequent rows down
# by one and dropping the previous last row (VeChain).

# Row 2
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '93.274.04'
$ws.Range("E2").Value = '  +1.77%  '

# Row 3
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.118.67'
$ws.Range("E3").Value = '  +0.29%  '

# Row 4
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.01%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '615.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.10'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.99%  '

# Row 8
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.391'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.68%  '

# Row 9
$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '

# Row 10
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.840'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +13.52%  '

# Row 11
$ws.Range("B11").Value = 'LidoStakedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D11").Value = '3.117.48'
$ws.Range("E11").Value = '  +0.37%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.198'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.36%  '

# Row 13
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.28%  '

# Row 14
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000246'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.38%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '93.001.39'
$ws.Range("E15").Value = '  +1.79%  '

# Row 16
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.82%  '

# Row 17
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '3.702.37'
$ws.Range("E17").Value = '  +0.51%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.106.10'
$ws.Range("E18").Value = '  +0.07%  '

# Row 19
$ws.Range("B19").Value = 'SuiNetwork'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.11%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '

# Row 21
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.93%  '

# Row 22
$ws.Range("B22").Value = 'PEPE'
$ws.Range("C22").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000202'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.67%  '

# Row 23
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '443.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '

# Row 24
$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.09%  '

# Row 25
$ws.Range("B25").Value = 'Aptos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.67%  '

# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.79%  '

# Row 27
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '86.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.92%  '

# Row 28
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '3.288.83'
$ws.Range("E28").Value = '  +0.05%  '

# Row 29
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '

# Row 30
$ws.Range("B30").Value = 'Cronos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.181'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.17%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.237'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.29%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.124'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -11.79%  '

# Row 33
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.94%  '

# Row 34
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.40%  '

# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.03'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.60%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.160'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.57%  '

# Row 37
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.03'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '

# Row 38
$ws.Range("B38").Value = 'MantraDAO'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.97%  '

# Row 39
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.02%  '

# Row 40
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.452'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.00%  '

# Row 41
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.43%  '

# Row 42
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '479.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.86%  '

# Row 43
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.99%  '

# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.69%  '

# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '159.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.93%  '

# Row 47
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.704'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.36%  '

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.20%  '

# Row 49
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.09%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.02%  '

# Row 51
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.21%  '

